$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> Data, add new sheet RF ---
$data = $wb.Worksheets.Item("Sheet1")
$data.Name = "Data"

$rf = $wb.Worksheets.Add($null, $data)
$rf.Name = "RF"

# --- Fill in the "Data" sheet ---

# Row 1: group numbers across G1:R1 (three columns per group)
$data.Range("G1").Value = 1
$data.Range("H1").Value = 1
$data.Range("I1").Value = 1
$data.Range("J1").Value = 2
$data.Range("K1").Value = 2
$data.Range("L1").Value = 2
$data.Range("M1").Value = 3
$data.Range("N1").Value = 3
$data.Range("O1").Value = 3
$data.Range("P1").Value = 4
$data.Range("Q1").Value = 4
$data.Range("R1").Value = 4

# Row 2: batch-correction method headers (B:F) and repeated clr headers (G:R)
$data.Range("B2").Value = "uncorrected"
$data.Range("C2").Value = "bmc"
$data.Range("D2").Value = "combat"
$data.Range("E2").Value = "limma"
$data.Range("F2").Value = "DCC"

$data.Range("G2").Value = "clr_pca1"
$data.Range("H2").Value = "clr_pcacounts"
$data.Range("I2").Value = "clr_pcaroundcounts"
$data.Range("J2").Value = "clr_pca1"
$data.Range("K2").Value = "clr_pcacounts"
$data.Range("L2").Value = "clr_pcaroundcounts"
$data.Range("M2").Value = "clr_pca1"
$data.Range("N2").Value = "clr_pcacounts"
$data.Range("O2").Value = "clr_pcaroundcounts"
$data.Range("P2").Value = "clr_pca1"
$data.Range("Q2").Value = "clr_pcacounts"
$data.Range("R2").Value = "clr_pcaroundcounts"

# Column A: dataset labels
$data.Range("A3").Value = "Gibbonsr_complete_otu"
$data.Range("A4").Value = "Thomasr_complete_otu"
$data.Range("A5").Value = "AGPr_complete_otu"
$data.Range("A6").Value = "AGPr_max_k5"
$data.Range("A7").Value = "AGPr_max_k6"
$data.Range("A8").Value = "AGPr_max_k7"

# Lone numeric values
$data.Range("J4").Value = 27
$data.Range("L4").Value = 27

# --- Column widths ---
# (values picked so the stored OOXML <col width=.../> lands as close as
# possible to the authored 33.83203125 / 16.83203125 / 15.6640625 widths,
# given this engine's column-width rounding model)
$data.Columns.Item(1).ColumnWidth = 33
$data.Range("G1:I1").EntireColumn.ColumnWidth = 16
$data.Columns.Item(11).ColumnWidth = 14.8

# --- Selection on Data sheet ---
$data.Activate()
$data.Range("L9").Select()

# --- Window position ---
$excel.ActiveWindow.Left = -30660
$excel.ActiveWindow.Top = 1420
